$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

# --- Correct a handful of previously-estimated figures now that final
#     conversions are in -------------------------------------------------
$ws.Cells.Item(189, 3).Value2 = 196
$ws.Cells.Item(189, 4).Value2 = 86.53
$ws.Cells.Item(196, 4).Value2 = 88.1
$ws.Cells.Item(200, 3).Value2 = 265
$ws.Cells.Item(200, 4).Value2 = 89.7

# --- Append the newly completed ridership conversions for 25-31 Jul 2016 -
$ws.Cells.Item(203, 1).Value2 = "Monday"
$ws.Cells.Item(203, 2).Value2 = "25 Jul 2016"
$ws.Cells.Item(203, 3).Value2 = 185
$ws.Cells.Item(203, 4).Value2 = 90.4
$ws.Cells.Item(203, 5).Value2 = 68.66

$ws.Cells.Item(204, 1).Value2 = "Tuesday"
$ws.Cells.Item(204, 2).Value2 = "26 Jul 2016"
$ws.Cells.Item(204, 3).Value2 = 236
$ws.Cells.Item(204, 4).Value2 = 91.98
$ws.Cells.Item(204, 5).Value2 = 68.81

$ws.Cells.Item(205, 1).Value2 = "Wednesday"
$ws.Cells.Item(205, 2).Value2 = "27 Jul 2016"
$ws.Cells.Item(205, 3).Value2 = 209
$ws.Cells.Item(205, 4).Value2 = 98.04
$ws.Cells.Item(205, 5).Value2 = 68.96

$ws.Cells.Item(206, 1).Value2 = "Thursday"
$ws.Cells.Item(206, 2).Value2 = "28 Jul 2016"
$ws.Cells.Item(206, 3).Value2 = 179
$ws.Cells.Item(206, 4).Value2 = 95.98
$ws.Cells.Item(206, 5).Value2 = 69.11

$ws.Cells.Item(207, 1).Value2 = "Friday"
$ws.Cells.Item(207, 2).Value2 = "29 Jul 2016"
$ws.Cells.Item(207, 3).Value2 = 206
$ws.Cells.Item(207, 4).Value2 = 92.34
$ws.Cells.Item(207, 5).Value2 = 69.26

$ws.Cells.Item(208, 1).Value2 = "Saturday"
$ws.Cells.Item(208, 2).Value2 = "30 Jul 2016"
$ws.Cells.Item(208, 3).Value2 = 122
$ws.Cells.Item(208, 4).Value2 = 39.73
$ws.Cells.Item(208, 5).Value2 = 69.41

$ws.Cells.Item(209, 1).Value2 = "Sunday"
$ws.Cells.Item(209, 2).Value2 = "31 Jul 2016"
$ws.Cells.Item(209, 3).Value2 = 100
$ws.Cells.Item(209, 4).Value2 = 32.98
$ws.Cells.Item(209, 5).Value2 = 69.56

# --- Extend the chart's series ranges to cover the new rows -------------
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$sc = $chart.SeriesCollection()

$s1 = $sc.Item(1)
$s1.Formula = '=SERIES("Ridership",Ridership!$B$2:$B$209,Ridership!$C$2:$C$209,1)'

$s2 = $sc.Item(2)
$s2.Formula = '=SERIES("Average",Ridership!$B$2:$B$209,Ridership!$D$2:$D$209,2)'

$s3 = $sc.Item(3)
$s3.Formula = '=SERIES("Pilot",Ridership!$B$2:$B$209,Ridership!$E$2:$E$209,3)'

# --- Slide the chart down so it still sits below the (now longer) table -
$co.Top = 3172.5
